$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking total (B11: 3 -> 5)
$ws.Range("B11").Value = 5

# Update total correct marks (B12: 51 -> 85)
$ws.Range("B12").Value = 85

# Update correct/total marks display (E12: "50/84" -> "85/140")
$ws.Range("E12").Value = "85/140"
